$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 130.955829
$ws.Range("H2").Value = 392.867487
$ws.Range("I2").Value = 0.5336535908353144
$ws.Range("J2").Value = 0.5336535908353144
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 16.535604
$ws.Range("N2").Value = 49.606812
$ws.Range("O2").Value = 0.2120453146491552
$ws.Range("P2").Value = 0.2120453146491552
$ws.Range("Q2").Value = 2165.433729835716
$ws.Range("R2").Value = 19488.90356852144
$ws.Range("S2").Value = 0.1131587435823258
$ws.Range("T2").Value = 0.1131587435823258

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 130.955829
$ws.Range("H3").Value = 392.867487
$ws.Range("I3").Value = 0.5336535908353144
$ws.Range("J3").Value = 0.5336535908353144
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 40.62063066666667
$ws.Range("N3").Value = 121.861892
$ws.Range("O3").Value = 0.5209011059384622
$ws.Range("P3").Value = 0.5209011059384622
$ws.Range("Q3").Value = 5319.508363456156
$ws.Range("R3").Value = 47875.57527110541
$ws.Range("S3").Value = 0.2779807456541468
$ws.Range("T3").Value = 0.2779807456541468

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 130.955829
$ws.Range("H4").Value = 392.867487
$ws.Range("I4").Value = 0.5336535908353144
$ws.Range("J4").Value = 0.5336535908353144
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 20.825229
$ws.Range("N4").Value = 62.475687
$ws.Range("O4").Value = 0.2670535794123827
$ws.Range("P4").Value = 0.2670535794123827
$ws.Range("Q4").Value = 2727.185127809841
$ws.Range("R4").Value = 24544.66615028857
$ws.Range("S4").Value = 0.1425141015988418
$ws.Range("T4").Value = 0.1425141015988418

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 66.39541
$ws.Range("H5").Value = 199.18623
$ws.Range("I5").Value = 0.2705656497465488
$ws.Range("J5").Value = 0.2705656497465488
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 16.535604
$ws.Range("N5").Value = 49.606812
$ws.Range("O5").Value = 0.2120453146491552
$ws.Range("P5").Value = 0.2120453146491552
$ws.Range("Q5").Value = 1097.88820717764
$ws.Range("R5").Value = 9880.993864598759
$ws.Range("S5").Value = 0.05737217833376007
$ws.Range("T5").Value = 0.05737217833376008

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 66.39541
$ws.Range("H6").Value = 199.18623
$ws.Range("I6").Value = 0.2705656497465488
$ws.Range("J6").Value = 0.2705656497465488
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 40.62063066666667
$ws.Range("N6").Value = 121.861892
$ws.Range("O6").Value = 0.5209011059384622
$ws.Range("P6").Value = 0.5209011059384622
$ws.Range("Q6").Value = 2697.023427571907
$ws.Range("R6").Value = 24273.21084814716
$ws.Range("S6").Value = 0.1409379461819359
$ws.Range("T6").Value = 0.1409379461819359

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 66.39541
$ws.Range("H7").Value = 199.18623
$ws.Range("I7").Value = 0.2705656497465488
$ws.Range("J7").Value = 0.2705656497465488
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 20.825229
$ws.Range("N7").Value = 62.475687
$ws.Range("O7").Value = 0.2670535794123827
$ws.Range("P7").Value = 0.2670535794123827
$ws.Range("Q7").Value = 1382.69961779889
$ws.Range("R7").Value = 12444.29656019001
$ws.Range("S7").Value = 0.0722555252308529
$ws.Range("T7").Value = 0.0722555252308529

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 48.043585
$ws.Range("H8").Value = 144.130755
$ws.Range("I8").Value = 0.1957807594181367
$ws.Range("J8").Value = 0.1957807594181367
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 16.535604
$ws.Range("N8").Value = 49.606812
$ws.Range("O8").Value = 0.2120453146491552
$ws.Range("P8").Value = 0.2120453146491552
$ws.Range("Q8").Value = 794.42969630034
$ws.Range("R8").Value = 7149.86726670306
$ws.Range("S8").Value = 0.04151439273306936
$ws.Range("T8").Value = 0.04151439273306936

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 48.043585
$ws.Range("H9").Value = 144.130755
$ws.Range("I9").Value = 0.1957807594181367
$ws.Range("J9").Value = 0.1957807594181367
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 40.62063066666667
$ws.Range("N9").Value = 121.861892
$ws.Range("O9").Value = 0.5209011059384622
$ws.Range("P9").Value = 0.5209011059384622
$ws.Range("Q9").Value = 1951.560722187607
$ws.Range("R9").Value = 17564.04649968846
$ws.Range("S9").Value = 0.1019824141023794
$ws.Range("T9").Value = 0.1019824141023794

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 48.043585
$ws.Range("H10").Value = 144.130755
$ws.Range("I10").Value = 0.1957807594181367
$ws.Range("J10").Value = 0.1957807594181367
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 20.825229
$ws.Range("N10").Value = 62.475687
$ws.Range("O10").Value = 0.2670535794123827
$ws.Range("P10").Value = 0.2670535794123827
$ws.Range("Q10").Value = 1000.518659605965
$ws.Range("R10").Value = 9004.667936453685
$ws.Range("S10").Value = 0.05228395258268796
$ws.Range("T10").Value = 0.05228395258268796
